$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1925.8182
$ws.Range("I40").Value = 1100
$ws.Range("J40").Value = 2235.5
$ws.Range("K40").Value = 1100
$ws.Range("L40").Value = 2235.5
$ws.Range("M40").Value = -925
$ws.Range("N40").Value = -2585.5
$ws.Range("H53").Value = 444.41666
$ws.Range("I53").Value = 462.75
$ws.Range("J53").Value = 435.25
$ws.Range("K53").Value = 462.75
$ws.Range("L53").Value = 435.25
$ws.Range("M53").Value = 174.25
$ws.Range("N53").Value = -1709.25
$ws.Range("H64").Value = 62504464
$ws.Range("I64").Value = 333342530
$ws.Range("J64").Value = 3375.4614
$ws.Range("K64").Value = 333342530
$ws.Range("L64").Value = 3375.4614
$ws.Range("M64").Value = -333342282
$ws.Range("N64").Value = -3871.4614
$ws.Range("H67").Value = 62504464
$ws.Range("I67").Value = 333342530
$ws.Range("J67").Value = 3375.4614
$ws.Range("K67").Value = 333342530
$ws.Range("L67").Value = 3375.4614
$ws.Range("M67").Value = -333341672
$ws.Range("N67").Value = -5091.4614
$ws.Range("H74").Value = 3231.3635
$ws.Range("I74").Value = 2535.5715
$ws.Range("J74").Value = 3418.6924
$ws.Range("K74").Value = 2535.5715
$ws.Range("L74").Value = 3418.6924
$ws.Range("M74").Value = -1599.5715
$ws.Range("N74").Value = -5290.6924
$ws.Range("H76").Value = 4089.862
$ws.Range("I76").Value = 3043.7368
$ws.Range("J76").Value = 6077.5
$ws.Range("K76").Value = 3043.7368
$ws.Range("L76").Value = 6077.5
$ws.Range("M76").Value = -2728.7368
$ws.Range("N76").Value = -6707.5
$ws.Range("H77").Value = 3231.3635
$ws.Range("I77").Value = 2535.5715
$ws.Range("J77").Value = 3418.6924
$ws.Range("K77").Value = 12677.8575
$ws.Range("L77").Value = 17093.462
$ws.Range("M77").Value = -7997.8575
$ws.Range("N77").Value = -26453.462
$ws.Range("H79").Value = 4089.862
$ws.Range("I79").Value = 3043.7368
$ws.Range("J79").Value = 6077.5
$ws.Range("K79").Value = 3043.7368
$ws.Range("L79").Value = 6077.5
$ws.Range("M79").Value = -1951.7368
$ws.Range("N79").Value = -8261.5
$ws.Range("H113").Value = 2569.75
$ws.Range("I113").Value = 2001.1818
$ws.Range("J113").Value = 3264.6667
$ws.Range("K113").Value = 2001.1818
$ws.Range("L113").Value = 3264.6667
$ws.Range("M113").Value = 1252.8182
$ws.Range("N113").Value = -9772.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8628.610000000001
$ws.Range("I32").Value = 6728.7754
$ws.Range("J32").Value = 24000
$ws.Range("K32").Value = 6728.7754
$ws.Range("L32").Value = 24000
$ws.Range("M32").Value = -6441.7754
$ws.Range("N32").Value = -24574
$ws.Range("H122").Value = 3856.6667
$ws.Range("I122").Value = 3762.5715
$ws.Range("J122").Value = 3988.4
$ws.Range("K122").Value = 11287.7145
$ws.Range("L122").Value = 11965.2
$ws.Range("M122").Value = -8837.7145
$ws.Range("N122").Value = -16865.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6386.0713
$ws.Range("I62").Value = 7880.5
$ws.Range("J62").Value = 2650
$ws.Range("K62").Value = 7880.5
$ws.Range("L62").Value = 2650
$ws.Range("M62").Value = -7256.5
$ws.Range("N62").Value = -3898
$ws.Range("H65").Value = 6386.0713
$ws.Range("I65").Value = 7880.5
$ws.Range("J65").Value = 2650
$ws.Range("K65").Value = 39402.5
$ws.Range("L65").Value = 13250
$ws.Range("M65").Value = -36282.5
$ws.Range("N65").Value = -19490
$ws.Range("H133").Value = 42000
$ws.Range("J133").Value = 42000
$ws.Range("L133").Value = 42000
$ws.Range("N133").Value = -47060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 958
$ws.Range("I5").Value = 753.6923
$ws.Range("J5").Value = 1162.3077
$ws.Range("K5").Value = 2261.0769
$ws.Range("L5").Value = 3486.9231
$ws.Range("M5").Value = -2149.0769
$ws.Range("N5").Value = -3710.9231
$ws.Range("H135").Value = 958
$ws.Range("I135").Value = 753.6923
$ws.Range("J135").Value = 1162.3077
$ws.Range("K135").Value = 6783.2307
$ws.Range("L135").Value = 10460.7693
$ws.Range("M135").Value = -4248.2307
$ws.Range("N135").Value = -15530.7693
$ws.Range("H137").Value = 3295264.2
$ws.Range("I137").Value = 56767.5
$ws.Range("K137").Value = 170302.5
$ws.Range("M137").Value = -165202.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3998.4285
$ws.Range("I70").Value = 3993.3333
$ws.Range("K70").Value = 3993.3333
$ws.Range("M70").Value = -3723.3333
$ws.Range("H73").Value = 3998.4285
$ws.Range("I73").Value = 3993.3333
$ws.Range("K73").Value = 3993.3333
$ws.Range("M73").Value = -3057.3333
$ws.Range("H80").Value = 3124.2
$ws.Range("I80").Value = 2205
$ws.Range("J80").Value = 4503
$ws.Range("K80").Value = 2205
$ws.Range("L80").Value = 4503
$ws.Range("M80").Value = -1207
$ws.Range("N80").Value = -6499
$ws.Range("H83").Value = 3124.2
$ws.Range("I83").Value = 2205
$ws.Range("J83").Value = 4503
$ws.Range("K83").Value = 11025
$ws.Range("L83").Value = 22515
$ws.Range("M83").Value = -6033
$ws.Range("N83").Value = -32499
$ws.Range("H132").Value = 4682.587
$ws.Range("I132").Value = 4391.8203
$ws.Range("J132").Value = 6302.5713
$ws.Range("K132").Value = 13175.4609
$ws.Range("L132").Value = 18907.7139
$ws.Range("M132").Value = -10645.4609
$ws.Range("N132").Value = -23967.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1017
$ws.Range("I93").Value = 737.0222
$ws.Range("J93").Value = 2066.9167
$ws.Range("K93").Value = 737.0222
$ws.Range("L93").Value = 2066.9167
$ws.Range("M93").Value = 510.9778
$ws.Range("N93").Value = -4562.9167
$ws.Range("H120").Value = 48000
$ws.Range("J120").Value = 48000
$ws.Range("L120").Value = 48000
$ws.Range("N120").Value = -57676

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9101.154
$ws.Range("I136").Value = 9567.5
$ws.Range("J136").Value = 3505
$ws.Range("K136").Value = 28702.5
$ws.Range("L136").Value = 10515
$ws.Range("M136").Value = -26152.5
$ws.Range("N136").Value = -15615
